$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Overwrite the data with the new "privacyType" / "post text" / "expected result" table ---

# Header row
$ws.Range("A1").Value2 = "privacyType"
$ws.Range("B1").Value2 = "post text"
$ws.Range("C1").Value2 = "expected result"

# Row 2
$ws.Range("A2").Value2 = "Only me"
$ws.Range("B2").Value2 = "Hello  only me post"
$ws.Range("C2").Value2 = "Shared with Only me"

# Row 3
$ws.Range("A3").Value2 = "Public"
$ws.Range("B3").Value2 = "Hello  Public post"
$ws.Range("C3").Value2 = "Shared with Public"

# --- Formatting ---

# A2:B3 already carry the pre-existing "JetBrains Mono / vertical-center" row style
# from the original table, so they are left untouched here (re-setting the same
# properties through the Font object would only churn the style table).

# New style for column C data cells: Consolas, size 7, color FFE8EAED
$ws.Range("C2:C3").Font.Name = "Consolas"
$ws.Range("C2:C3").Font.Size = 7
$ws.Range("C2:C3").Font.Color = 15592168

# --- Column widths (re-sized now that the columns hold different content) ---
$ws.Columns.Item(1).ColumnWidth = 9.8
$ws.Columns.Item(2).ColumnWidth = 21.6
$ws.Columns.Item(3).ColumnWidth = 14.8

# --- Selection ---
$ws.Range("C8").Select() | Out-Null
